$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: title and link update
$ws.Range("D6").Value = "[프로그래머스 - Python] 2020 KAKAO BLIND RECRUITMENT > 문자열 압축"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/%ED%94%84%EB%A1%9C%EA%B7%B8%EB%9E%98%EB%A8%B8%EC%8A%A4-Python-2020-KAKAO-BLIND-RECRUITMENT-%EB%AC%B8%EC%9E%90%EC%97%B4-%EC%95%95%EC%B6%95"

# Row 26: title update only
$ws.Range("D26").Value = "ai plus(est soft)"

# Row 51: title and link update
$ws.Range("D51").Value = "[pyqt5] 야구 중계에 스트라이크 존이 제공되지 않을 때 시청자를 위한 가이드 앱"
$ws.Range("E51").Value = "https://bskyvision.com/1251"
